# Fix unit bug in psc detection scripts:
# descriptions incorrectly said "seconds" when the underlying fields
# (psc_start_ms, psc_risetime_ms, psc_decay_ms) are actually in milliseconds.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C6").Value = "the vector denoting the starting time of PSC in miliseconds"
$ws.Range("C8").Value = "rise time of the PSC in milisecods"
$ws.Range("C9").Value = "decay time of the PSC in miliseconds"

$ws.Range("G11").Select()
